# Update workbook with the latest daily "Tasas de captación" rows.
# (Actualización desde MV -datos-)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dd-mm-yyyy text labels. Pre-format the target cells as
# text so Excel does not reinterpret the strings as date serials.
$ws.Range("A174:A179").NumberFormat = "@"

$ws.Range("A174").Value = "08-09-2021"
$ws.Range("B174").Value = 0.13
$ws.Range("C174").Value = 0.21
$ws.Range("D174").Value = 0
$ws.Range("E174").Value = 0.26
$ws.Range("F174").Value = 0.3

$ws.Range("A175").Value = "09-09-2021"
$ws.Range("B175").Value = 0.13
$ws.Range("C175").Value = 0.23
$ws.Range("D175").Value = 0.01
$ws.Range("E175").Value = 0.23
$ws.Range("F175").Value = 0.26

$ws.Range("A176").Value = "10-09-2021"
$ws.Range("B176").Value = 0.14
$ws.Range("C176").Value = 0.24
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 0.21
$ws.Range("F176").Value = 0.26

$ws.Range("A177").Value = "13-09-2021"
$ws.Range("B177").Value = 0.13
$ws.Range("C177").Value = 0.16
$ws.Range("D177").Value = 0.01
$ws.Range("E177").Value = 0.32
$ws.Range("F177").Value = 0.23

$ws.Range("A178").Value = "14-09-2021"
$ws.Range("B178").Value = 0.14
$ws.Range("C178").Value = 0.23
$ws.Range("D178").Value = 0.01
$ws.Range("E178").Value = 0.31
$ws.Range("F178").Value = 0.28

$ws.Range("A179").Value = "15-09-2021"
$ws.Range("B179").Value = 0.14
$ws.Range("C179").Value = 0.22
$ws.Range("D179").Value = 0.01
$ws.Range("E179").Value = 0.14
$ws.Range("F179").Value = 0.29

# Restore the plain default style on column A so the new date-label cells
# match the rest of the sheet (no explicit cell style reference).
$ws.Range("A174:A179").Style = "Normal"
